# Update computed line-flow results (pl_mw.xlsx, "case with 380 kV done")
# Only the numeric data cells in columns B..O (excluding the always-zero
# columns G/K/L/N and the index column A) for rows 2-25 change value;
# everything else in the sheet is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9064210071257435
$ws.Range("C2").Value = 0.2203067207800018
$ws.Range("D2").Value = 0.2023424355200234
$ws.Range("E2").Value = 0.1547909353960861
$ws.Range("F2").Value = 1.17751814677009
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("I2").Value = 0.5231711738452223
$ws.Range("J2").Value = 0.1596138729009482
$ws.Range("M2").Value = 0.3703977477924028
$ws.Range("O2").Value = 2.689013073750971
$ws.Range("B3").Value = 0.8031573600239881
$ws.Range("C3").Value = 0.1929127848293888
$ws.Range("D3").Value = 0.200038364680232
$ws.Range("E3").Value = 0.1548333156997757
$ws.Range("F3").Value = 1.183726354629734
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("I3").Value = 0.5328619925022533
$ws.Range("J3").Value = 0.1609187147277886
$ws.Range("M3").Value = 0.343004612708306
$ws.Range("O3").Value = 2.711440572946984
$ws.Range("B4").Value = 0.7396247032714314
$ws.Range("C4").Value = 0.1760380690245142
$ws.Range("D4").Value = 0.1986851779473184
$ws.Range("E4").Value = 0.1549289865395451
$ws.Range("F4").Value = 1.18835840074771
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("I4").Value = 0.5392324207284709
$ws.Range("J4").Value = 0.1618099164969031
$ws.Range("M4").Value = 0.3262333129422217
$ws.Range("O4").Value = 2.727302318494409
$ws.Range("B5").Value = 0.7137039781734131
$ws.Range("C5").Value = 0.1691481309999858
$ws.Range("D5").Value = 0.1981492872878832
$ws.Range("E5").Value = 0.154985507000525
$ws.Range("F5").Value = 1.190452051768879
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("I5").Value = 0.5419339037214268
$ws.Range("J5").Value = 0.1621957163346259
$ws.Range("M5").Value = 0.3194114119667333
$ws.Range("O5").Value = 2.734291136818598
$ws.Range("B6").Value = 0.7093980583691462
$ws.Range("C6").Value = 0.1680032675435541
$ws.Range("D6").Value = 0.198061244042286
$ws.Range("E6").Value = 0.1549959517687185
$ws.Range("F6").Value = 1.190812140644859
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("I6").Value = 0.5423888482487325
$ws.Range("J6").Value = 0.1622611445010484
$ws.Range("M6").Value = 0.3182794103298789
$ws.Range("O6").Value = 2.735483311316628
$ws.Range("B7").Value = 0.7392752492501131
$ws.Range("C7").Value = 0.1759452024378163
$ws.Range("D7").Value = 0.198677887702722
$ws.Range("E7").Value = 0.1549296777738753
$ws.Range("F7").Value = 1.188385802400113
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("I7").Value = 0.5392684270312671
$ws.Range("J7").Value = 0.1618150279212394
$ws.Range("M7").Value = 0.3261412590663326
$ws.Range("O7").Value = 2.727394447409537
$ws.Range("B8").Value = 0.870843282234091
$ws.Range("C8").Value = 0.2108728857933784
$ws.Range("D8").Value = 0.2015352660586842
$ws.Range("E8").Value = 0.1547911040532526
$ws.Range("F8").Value = 1.179488394169503
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("I8").Value = 0.5264252384303365
$ws.Range("J8").Value = 0.1600450923155066
$ws.Range("M8").Value = 0.3609428344195962
$ws.Range("O8").Value = 2.696311561930116
$ws.Range("B9").Value = 1.127770804551631
$ws.Range("C9").Value = 0.2789171793601213
$ws.Range("D9").Value = 0.2076239239028723
$ws.Range("E9").Value = 0.1550712818679756
$ws.Range("F9").Value = 1.168557834434139
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("I9").Value = 0.50458221266123
$ws.Range("J9").Value = 0.1572890219108487
$ws.Range("M9").Value = 0.4295561425678116
$ws.Range("O9").Value = 2.651988839445721
$ws.Range("B10").Value = 1.315819253883262
$ws.Range("C10").Value = 0.3286211394548388
$ws.Range("D10").Value = 0.2123900162200414
$ws.Range("E10").Value = 0.1556128750754873
$ws.Range("F10").Value = 1.164514829431837
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("I10").Value = 0.4905825250121119
$ws.Range("J10").Value = 0.1557006066866329
$ws.Range("M10").Value = 0.480175673445153
$ws.Range("O10").Value = 2.629618584401015
$ws.Range("B11").Value = 1.401200570918093
$ws.Range("C11").Value = 0.3511674419906967
$ws.Range("D11").Value = 0.2146211864236705
$ws.Range("E11").Value = 0.1559320291362667
$ws.Range("F11").Value = 1.163544680145577
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("I11").Value = 0.4846607805398371
$ws.Range("J11").Value = 0.1550729349000797
$ws.Range("M11").Value = 0.503246512221196
$ws.Range("O11").Value = 2.62166720806627
$ws.Range("B12").Value = 1.433507491900457
$ws.Range("C12").Value = 0.3596955626244949
$ws.Range("D12").Value = 0.2154750747441199
$ws.Range("E12").Value = 0.1560633346872997
$ws.Range("F12").Value = 1.163302512095342
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("I12").Value = 0.4824828185386494
$ws.Range("J12").Value = 0.1548489124218762
$ws.Range("M12").Value = 0.5119887883586784
$ws.Range("O12").Value = 2.618977114009084
$ws.Range("B13").Value = 1.426550763198122
$ws.Range("C13").Value = 0.3578593163812229
$ws.Range("D13").Value = 0.2152907758441245
$ws.Range("E13").Value = 0.1560345912419976
$ws.Range("F13").Value = 1.163349094915247
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("I13").Value = 0.4829490105075251
$ws.Range("J13").Value = 0.1548965517444429
$ws.Range("M13").Value = 0.5101057295947697
$ws.Range("O13").Value = 2.619542185470692
$ws.Range("B14").Value = 1.403858994155144
$ws.Range("C14").Value = 0.3518692519646152
$ws.Range("D14").Value = 0.2146912565273027
$ws.Range("E14").Value = 0.1559426224095866
$ws.Range("F14").Value = 1.163522246139166
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("I14").Value = 0.4844803048381934
$ws.Range("J14").Value = 0.155054230514672
$ws.Range("M14").Value = 0.5039656298167756
$ws.Range("O14").Value = 2.621439455147708
$ws.Range("B15").Value = 1.389956307120883
$ws.Range("C15").Value = 0.3481988918565548
$ws.Range("D15").Value = 0.2143252025441882
$ws.Range("E15").Value = 0.1558876490892693
$ws.Range("F15").Value = 1.163644618946677
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("I15").Value = 0.4854266699989473
$ws.Range("J15").Value = 0.1551525931350746
$ws.Range("M15").Value = 0.5002053902025807
$ws.Range("O15").Value = 2.622643409599561
$ws.Range("B16").Value = 1.310235960142506
$ws.Range("C16").Value = 0.3271463551298268
$ws.Range("D16").Value = 0.2122454664080635
$ws.Range("E16").Value = 0.1555934806908503
$ws.Range("F16").Value = 1.16459573873756
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("I16").Value = 0.4909785347967848
$ws.Range("J16").Value = 0.1557435376212943
$ws.Range("M16").Value = 0.4786687825430747
$ws.Range("O16").Value = 2.630183080739329
$ws.Range("B17").Value = 1.261287173422033
$ws.Range("C17").Value = 0.3142145341875278
$ws.Range("D17").Value = 0.2109857116946188
$ws.Range("E17").Value = 0.1554316467074344
$ws.Range("F17").Value = 1.165401971881948
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("I17").Value = 0.4944990177940269
$ws.Range("J17").Value = 0.1561303819387767
$ws.Range("M17").Value = 0.4654676635006467
$ws.Range("O17").Value = 2.635378984821472
$ws.Range("B18").Value = 1.233117897545696
$ws.Range("C18").Value = 0.306770468213756
$ws.Range("D18").Value = 0.2102670731579366
$ws.Range("E18").Value = 0.1553454155445166
$ws.Range("F18").Value = 1.165947479166682
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("I18").Value = 0.4965659463327619
$ws.Range("J18").Value = 0.1563618166412013
$ws.Range("M18").Value = 0.4578788741522004
$ws.Range("O18").Value = 2.638576897762732
$ws.Range("B19").Value = 1.22357770163768
$ws.Range("C19").Value = 0.3042490138643359
$ws.Range("D19").Value = 0.2100247767987753
$ws.Range("E19").Value = 0.1553173964162902
$ws.Range("F19").Value = 1.16614621757634
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("I19").Value = 0.4972729862771068
$ws.Range("J19").Value = 0.1564417100676749
$ws.Range("M19").Value = 0.4553101659446668
$ws.Range("O19").Value = 2.639695585323096
$ws.Range("B20").Value = 1.266499443855196
$ws.Range("C20").Value = 0.3155917752675066
$ws.Range("D20").Value = 0.2111192006029938
$ws.Range("E20").Value = 0.1554481652332989
$ws.Range("F20").Value = 1.165307680974223
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("I20").Value = 0.4941199037559514
$ws.Range("J20").Value = 0.1560882771490739
$ws.Range("M20").Value = 0.4668725193657082
$ws.Range("O20").Value = 2.634804196268107
$ws.Range("B21").Value = 1.4105248123808
$ws.Range("C21").Value = 0.353628945080402
$ws.Range("D21").Value = 0.214867106415582
$ws.Range("E21").Value = 0.1559693524413426
$ws.Range("F21").Value = 1.16346798745198
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("I21").Value = 0.4840287749567302
$ws.Range("J21").Value = 0.1550075454960513
$ws.Range("M21").Value = 0.5057689706973463
$ws.Range("O21").Value = 2.620873463658285
$ws.Range("B22").Value = 1.50450645199362
$ws.Range("C22").Value = 0.3784318444796213
$ws.Range("D22").Value = 0.2173689610878853
$ws.Range("E22").Value = 0.1563708743886174
$ws.Range("F22").Value = 1.162995550838531
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("I22").Value = 0.4778095485602769
$ws.Range("J22").Value = 0.1543808718294741
$ws.Range("M22").Value = 0.5312239465486925
$ws.Range("O22").Value = 2.613639937952343
$ws.Range("B23").Value = 1.454360718268902
$ws.Range("C23").Value = 0.3651993899205195
$ws.Range("D23").Value = 0.2160289056247819
$ws.Range("E23").Value = 0.1561510081648727
$ws.Range("F23").Value = 1.163180833325839
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("I23").Value = 0.4810943960192375
$ws.Range("J23").Value = 0.1547080464399926
$ws.Range("M23").Value = 0.5176351907092993
$ws.Range("O23").Value = 2.617329088694362
$ws.Range("B24").Value = 1.264143063377105
$ws.Range("C24").Value = 0.31496915379347
$ws.Range("D24").Value = 0.2110588327841185
$ws.Range("E24").Value = 0.1554406759962319
$ws.Range("F24").Value = 1.165350054506533
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("I24").Value = 0.4942911674227943
$ws.Range("J24").Value = 0.1561072845896057
$ws.Range("M24").Value = 0.4662373817858665
$ws.Range("O24").Value = 2.635063401864585
$ws.Range("B25").Value = 1.058386627012453
$ws.Range("C25").Value = 0.2605589639082666
$ws.Range("D25").Value = 0.20592514533746
$ws.Range("E25").Value = 0.1549364739725156
$ws.Range("F25").Value = 1.170815612382199
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("I25").Value = 0.5101323877939912
$ws.Range("J25").Value = 0.1579580159244891
$ws.Range("M25").Value = 0.4109566133679579
$ws.Range("O25").Value = 2.662193166503272
